# Edit powerpoint for diagrams and update png diagrams
#
# Renames the placeholder term "address" -> "giatros" in the two activity-
# diagram callouts on slide 1:
#   1. "[command commits address book]"                       -> "[command commits giatros book]"
#   2. "Purge redundant states and then save address book to
#       addressBookStateList"                                 -> "...save giatros book to giatrosBookStateList"
#
# (Only the slide content is touched here; the surrounding date placeholder
# caches / revisionInfo.xml bookkeeping file seen in the original commit are
# app-managed, resave-time artifacts that aren't part of the PowerPoint
# object model, so they aren't reproduced by this script.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape "TextBox 47": "[command commits address book]" ---------------
$sh1 = $s.Shapes.Item("TextBox 47")
$tr1 = $sh1.TextFrame.TextRange
# Replace just the word "address" (chars 18-24) with "giatros", leaving the
# surrounding runs ("[" and the rest of the sentence) untouched.
$tr1.Characters(18, 7).Text = "giatros"

# --- Shape "Rounded Rectangle 50": "Purge redundant states..." ----------
$sh2 = $s.Shapes.Item("Rounded Rectangle 50")
$tr2 = $sh2.TextFrame.TextRange
# Replace the later occurrence first so the earlier offset stays valid.
$tr2.Characters(54, 20).Text = "giatrosBookStateList"
$tr2.Characters(38, 7).Text = "giatros"
